$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date value from 45192 to 45202 for all data rows (2-199)
for ($i = 2; $i -le 199; $i++) {
    $ws.Cells.Item($i, 3).Value = 45202
}

# Row 199 gains an explicit custom row height (matches rows around it)
$ws.Rows.Item(199).RowHeight = 15

# Append new row 200 with the new logging entry
$ws.Cells.Item(200, 1).Value = "A 46301-2023"

$ws.Cells.Item(200, 2).Value = 45197
$ws.Cells.Item(200, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(200, 3).Value = 45202
$ws.Cells.Item(200, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(200, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(200, 5).Value = "HAPARANDA"

$ws.Cells.Item(200, 7).Value = 0.6
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(200, 9).Value = 0
$ws.Cells.Item(200, 10).Value = 0
$ws.Cells.Item(200, 11).Value = 0
$ws.Cells.Item(200, 12).Value = 0
$ws.Cells.Item(200, 13).Value = 0
$ws.Cells.Item(200, 14).Value = 0
$ws.Cells.Item(200, 15).Value = 0
$ws.Cells.Item(200, 16).Value = 0
$ws.Cells.Item(200, 17).Value = 0

$ws.Cells.Item(200, 18).Value = ""
$ws.Cells.Item(200, 18).WrapText = $true
